# TC02_Canine_Filter_SamplePatho-BCellLymphoma.xlsx
#
# The "Cases" query stored in cell B2 of the "startup" sheet is updated to
# drop the OPTIONAL MATCH cohort lookup and the trailing `Cohort` column
# from its RETURN clause (the last RETURN line loses its trailing comma
# too). The "Samples" (B3) and "Files" (B4) query text is unchanged.
#
# Finally, the active selection on the sheet moves back to B2 (it had
# drifted to B4 in the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["B Cell Lymphoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

$ws.Range("B2").Select() | Out-Null
